$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" column (col 7 / G) values, replacing the old Strike# derived
# values with newly regenerated std/mean based s_vals.
$kValues = @{
    2 = 0
    3 = 0
    4 = 3
    5 = 2
    6 = 3
    7 = 2
    8 = 3
    9 = 3
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    22 = 2
    23 = 1
    24 = 2
    25 = 1
    26 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 2
    35 = 2
    36 = 0
    37 = 1
    38 = 0
    39 = 2
    40 = 2
    41 = 1
    42 = 2
    43 = 2
    44 = 2
    45 = 2
    46 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 1
    52 = 1
    53 = 2
    54 = 3
    55 = 2
    56 = 1
    57 = 0
    58 = 0
    59 = 2
    60 = 1
    61 = 2
    62 = 2
    63 = 1
    64 = 0
    65 = 2
    66 = 2
    67 = 1
    68 = 2
    69 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}

